$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.846.77"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.598.43"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.46%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -4.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.245"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0610"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0785"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "1.821.37"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").Value = "1.600.45"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.507"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.47%  "
$ws.Range("D16").Value = "25.844.31"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("E18").Value = "  -4.15%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "189.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.128"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0460"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("E33").Value = "  -5.03%  "
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("D36").Value = "1.101.49"
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("E37").Value = "  -2.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.795"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.494"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "95.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("D42").Value = "1.733.83"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.740"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.87%  "
$ws.Range("D45").Value = "0.0₆0101"
$ws.Range("E45").Value = "  -11.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "52.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.89%  "
